# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.854.47'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '2.234.29'
$ws.Range('E3').Value = '  -2.42%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.649'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '230.04'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '62.94'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.37%  '
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.448'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0958'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.71%  '
$ws.Range('E11').Value = '  -1.83%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '26.48'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +9.05%  '
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').Value = '2.564.16'
$ws.Range('E14').Value = '  -2.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.38'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.10'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.822'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('D18').Value = '2.233.37'
$ws.Range('E18').Value = '  -2.69%  '
$ws.Range('D19').Value = '43.613.85'
$ws.Range('E19').Value = '  -0.85%  '
$ws.Range('D20').Value = '0.0₃0984'
$ws.Range('E20').Value = '  +4.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.59'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '248.18'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.40%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.40'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -6.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.38'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +22.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.23'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.82'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.37'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.53%  '
$ws.Range('E30').Value = '  +0.92%  '
$ws.Range('E31').Value = '  -1.77%  '
$ws.Range('E32').Value = '  -3.79%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.126'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.09%  '
$ws.Range('E34').Value = '  +5.40%  '
$ws.Range('E35').Value = '  -1.75%  '
$ws.Range('E36').Value = '  -3.42%  '
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('E38').Value = '  -2.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.27'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0258'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.95%  '
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('E42').Value = '  -2.40%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.24'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.04'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '96.78'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.15%  '
$ws.Range('E46').Value = '  -2.46%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.18'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.35%  '
$ws.Range('B48').Value = 'FTXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.35'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.31'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.84%  '
$ws.Range('D50').Value = '1.428.74'
$ws.Range('E50').Value = '  -3.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.76'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.42%  '
